$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 2 (B2 and C2) without shifting other rows
$ws.Range("B2:C2").ClearContents()

# Update B4 value from 6.62 to 1
$ws.Range("B4").Value = 1
